$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are numeric-looking strings need to be forced to
# stay text (matching the source inline-string cells) - otherwise Excel
# auto-converts them to numbers. Apply a temporary Text number format,
# assign the value, then restore the default style so no visible
# formatting change is left behind.
$textForceCells = @(
    @{Cell='D5'; Value='216.73'},
    @{Cell='D6'; Value='0.521'},
    @{Cell='D11'; Value='0.0846'},
    @{Cell='D15'; Value='0.541'},
    @{Cell='D16'; Value='66.67'},
    @{Cell='D19'; Value='217.06'},
    @{Cell='D21'; Value='6.91'},
    @{Cell='D22'; Value='4.42'},
    @{Cell='D25'; Value='146.25'},
    @{Cell='D29'; Value='15.68'},
    @{Cell='D37'; Value='0.0176'},
    @{Cell='D38'; Value='0.857'},
    @{Cell='D41'; Value='0.810'},
    @{Cell='D43'; Value='5.28'},
    @{Cell='D45'; Value='61.66'},
    @{Cell='D46'; Value='91.34'},
    @{Cell='D51'; Value='0.0960'}
)

foreach ($item in $textForceCells) {
    $rng = $ws.Range($item.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $item.Value
    $rng.Style = "Normal"
}

# Remaining cells already contain non-numeric-looking text (extra dots,
# percent signs, subscript digits, etc.) so a plain assignment keeps them
# stored as text without any style manipulation.
$ws.Range('D2').Value = '27.095.60'
$ws.Range('E2').Value = '  +0.99%  '
$ws.Range('D3').Value = '1.638.39'
$ws.Range('E3').Value = '  +0.06%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('E5').Value = '  +0.01%  '
$ws.Range('E6').Value = '  +2.08%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  -0.13%  '
$ws.Range('E9').Value = '  +0.13%  '
$ws.Range('E10').Value = '  +0.68%  '
$ws.Range('E11').Value = '  -0.02%  '
$ws.Range('D12').Value = '1.866.11'
$ws.Range('D13').Value = '1.633.12'
$ws.Range('E13').Value = '  -0.28%  '
$ws.Range('E14').Value = '  +0.32%  '
$ws.Range('E15').Value = '  +2.30%  '
$ws.Range('E16').Value = '  -0.77%  '
$ws.Range('D17').Value = '27.104.38'
$ws.Range('E17').Value = '  +1.07%  '
$ws.Range('E18').Value = '  +1.34%  '
$ws.Range('E19').Value = '  -0.45%  '
$ws.Range('E20').Value = '  +0.10%  '
$ws.Range('E21').Value = '  +1.57%  '
$ws.Range('E22').Value = '  +0.67%  '
$ws.Range('E23').Value = '  +3.19%  '
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('E25').Value = '  -0.45%  '
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('E27').Value = '  +1.41%  '
$ws.Range('E28').Value = '  +0.13%  '
$ws.Range('E29').Value = '  -0.46%  '
$ws.Range('E30').Value = '  +1.26%  '
$ws.Range('E31').Value = '  +0.21%  '
$ws.Range('E32').Value = '  +1.70%  '
$ws.Range('E33').Value = '  +0.78%  '
$ws.Range('D34').Value = '1.308.17'
$ws.Range('E34').Value = '  +3.33%  '
$ws.Range('E35').Value = '  +0.22%  '
$ws.Range('E36').Value = '  +1.18%  '
$ws.Range('E37').Value = '  -0.62%  '
$ws.Range('E38').Value = '  +3.01%  '
$ws.Range('E39').Value = '  +1.54%  '
$ws.Range('E40').Value = '  -0.03%  '
$ws.Range('E41').Value = '  +0.43%  '
$ws.Range('E42').Value = '  +5.33%  '
$ws.Range('E43').Value = '  -1.74%  '
$ws.Range('D44').Value = '1.777.02'
$ws.Range('E44').Value = '  -0.14%  '
$ws.Range('E45').Value = '  -0.08%  '
$ws.Range('E46').Value = '  -0.29%  '
$ws.Range('E47').Value = '  +1.00%  '
$ws.Range('D48').Value = '0.0₆0107'
$ws.Range('E48').Value = '  +1.25%  '
$ws.Range('E49').Value = '  -0.16%  '
$ws.Range('E50').Value = '  +0.26%  '
$ws.Range('E51').Value = '  +0.06%  '
